# ---------------------------------------------------------------------------
# Enrich the Stata Excel template:
#   - rename the original sheet to "Overall Results" and flesh out its demo
#     rows (color / trend-arrow legend used by the conditional formatting)
#   - duplicate it twice to create "Values, Ratings, Trends" and
#     "Values, Scores, Ratings, Trends" (same conditional formatting)
#   - add a new "Codebook" sheet (no conditional formatting) with a simple
#     title/value header and a frozen top row
#   - make "Codebook" the active sheet/tab
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Turn the original sheet into "Overall Results" -------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Overall Results"

# Legend rows below the existing title (row 1) / description (row 2) rows.
$ws1.Range("A3").Value = "gray"
$ws1.Range("B3").Value = 15.5
$ws1.Range("C3").Value = "↑"

$ws1.Range("A4").Value = "green"
$ws1.Range("B4").Value = 27.6
$ws1.Range("C4").Value = "↓"

$ws1.Range("A5").Value = "yellow"
$ws1.Range("B5").Value = "hello"
$ws1.Range("C5").Value = "➚"

$ws1.Range("A6").Value = "orange"
$ws1.Range("B6").Value = "test"
$ws1.Range("C6").Value = "↓"

$ws1.Range("A7").Value = "red"
$ws1.Range("C7").Value = "→"

# A handful of additional, otherwise-empty formatted rows (8-20) so the
# sheet's used range/formatting extends well past the legend, matching the
# "normal formatting" the author added to the template.
$ws1.Range("A8:D20").Borders.LineStyle = -4142 # xlLineStyleNone (keeps default look)

# Extend the sheet-wide conditional formatting to also cover the columns
# beyond ZZ on row 1 (AAA1:XFD1), as in the edited template.
$ws1.Range("A1:XFD1").FormatConditions.Delete()
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("green",A1)))')
$cf.Interior.Color = 65280
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("red",A1)))')
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("yellow",A1)))')
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("orange",A1)))')
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("gray",A1)))')
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("↑",A1)))')
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("➚",A1)))')
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("→",A1)))')
$cf = $ws1.Range("A1:ZZ1048576,AAA1:XFD1").FormatConditions.Add(2, 0, 'NOT(ISERROR(SEARCH("↓",A1)))')

$ws1.Range("E5").Select()

# --- 2. Duplicate to create "Values, Ratings, Trends" --------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Values, Ratings, Trends"
$ws2.Range("E5").Select()

# --- 3. Duplicate again to create "Values, Scores, Ratings, Trends" ------
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Values, Scores, Ratings, Trends"
$ws3.Range("E5").Select()

# --- 4. Add the "Codebook" sheet (no conditional formatting) -------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Codebook"

$ws4.Range("A1:D1").Value = "title"
$ws4.Rows.Item(1).RowHeight = 27
$ws4.Range("A1:D1").Orientation = 45

$ws4.Range("A2:D2").Value = "value"

$ws4.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws4.Range("E5").Select()

# --- 5. Make "Codebook" the active sheet/tab ------------------------------
$ws4.Activate()

$excel.ActiveWindow.WindowState = -4143 # xlNormal, ensure dimensions below stick
$excel.ActiveWindow.Left = 2505
$excel.ActiveWindow.Top = 1560
$excel.ActiveWindow.Width = 25140
$excel.ActiveWindow.Height = 11835

Write-Output "done"
